# "filling marker info into status 8"
# Fill the marker_1 column (J) with "NAT" for the rows that were still missing it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J5").Value = "NAT"
$ws.Range("J6").Value = "NAT"
$ws.Range("J7").Value = "NAT"
$ws.Range("J8").Value = "NAT"

# Leave the selection where the editor ended up clicking.
$ws.Range("Q11").Select()
